$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Find-ParaIndexByText($needle) {
    $idx = 0
    $found = -1
    foreach ($p in $d.Paragraphs) {
        $idx = $idx + 1
        if ($p.Range.Text -like $needle) {
            $found = $idx
            break
        }
    }
    return $found
}

# ---------------------------------------------------------------------------
# 1) Paragraph "3.3 Các quy trình nghiệp vụ cần giải quyết": drop the hidden
#    "_GoBack" bookmark that wraps the start of the paragraph.
# ---------------------------------------------------------------------------
$idx = Find-ParaIndexByText("*quy trình nghiệp vụ cần giải quyết*")
$p = $d.Paragraphs($idx)
$xml = "<w:p $wns>" +
         "<w:r><w:rPr><w:lang w:val=`"vi-VN`"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space=`"preserve`">3.3 </w:t></w:r>" +
         "<w:r><w:t>Các quy trình nghiệp vụ cần giải quyết</w:t></w:r>" +
       "</w:p>"
[void]$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 2) The three empty paragraphs right after it become five paragraphs that
#    list the business processes to cover.
# ---------------------------------------------------------------------------
$idx = Find-ParaIndexByText("*quy trình nghiệp vụ cần giải quyết*")
$pEmpty1 = $d.Paragraphs($idx + 1)
$xml1 = "<w:p $wns><w:r><w:t>Quy trình xử lý đăng nhập</w:t></w:r></w:p>"
[void]$pEmpty1.Range.InsertXML($xml1)

$idx = Find-ParaIndexByText("*Quy trình xử lý đăng nhập*")
$pEmpty2 = $d.Paragraphs($idx + 1)
$xml2 = "<w:p $wns><w:r><w:t>Quy trình xử lý đăng ký</w:t></w:r></w:p>"
[void]$pEmpty2.Range.InsertXML($xml2)

$idx = Find-ParaIndexByText("*Quy trình xử lý đăng ký*")
$pEmpty3 = $d.Paragraphs($idx + 1)
$xml3 = "<w:p $wns><w:r><w:t>Quy trình xử lý bán vé</w:t></w:r></w:p>" +
        "<w:p $wns><w:r><w:t>Quy trình xử lý hủy vé</w:t></w:r></w:p>" +
        "<w:p $wns><w:r><w:t>Q</w:t></w:r><w:r><w:t>uy trình xử lý thống kê số vé</w:t></w:r></w:p>"
[void]$pEmpty3.Range.InsertXML($xml3)

# ---------------------------------------------------------------------------
# 3) "Thao tác nghiệp vụ ( theo phiên bản thương mại)." gets split into three
#    runs with proofErr grammar-check markers bracketing "( theo".
# ---------------------------------------------------------------------------
$idx = Find-ParaIndexByText("*nghiệp vụ ( theo phiên bản thương mại*")
$p = $d.Paragraphs($idx)
$pPr = "<w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"3`"/></w:numPr></w:pPr>"
$xml = "<w:p $wns>" + $pPr +
         "<w:r><w:t xml:space=`"preserve`">Thao tác nghiệp vụ </w:t></w:r>" +
         "<w:proofErr w:type=`"gramStart`"/>" +
         "<w:r><w:t>( theo</w:t></w:r>" +
         "<w:proofErr w:type=`"gramEnd`"/>" +
         "<w:r><w:t xml:space=`"preserve`"> phiên bản thương mại).</w:t></w:r>" +
       "</w:p>"
[void]$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 4) "Hỗ trợ giao diện bằng các Component tùy chọn (Devpress, Dotnetbar,
#    Syncfusion..)" gets split the same way around "Syncfusion..".
# ---------------------------------------------------------------------------
$idx = Find-ParaIndexByText("*Dotnetbar, Syncfusion*")
$p = $d.Paragraphs($idx)
$pPr = "<w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"3`"/></w:numPr></w:pPr>"
$xml = "<w:p $wns>" + $pPr +
         "<w:r><w:t xml:space=`"preserve`">Hỗ trợ giao diện bằng các Component tùy chọn (Devpress, Dotnetbar, </w:t></w:r>" +
         "<w:proofErr w:type=`"gramStart`"/>" +
         "<w:r><w:t>Syncfusion..</w:t></w:r>" +
         "<w:proofErr w:type=`"gramEnd`"/>" +
         "<w:r><w:t>)</w:t></w:r>" +
       "</w:p>"
[void]$p.Range.InsertXML($xml)
